# Add the new "2022-Q3" quarterly sheet and update the "总计" (total)
# summary sheet to include it as the first (most recent) row.
#
# NOTE: this runtime's Range.Value *getter* returns a bogus placeholder
# when read back into a variable, so any value that needs to survive a
# read must instead be fetched through .Formula (which round-trips
# correctly for both text and numbers) - or, where the value is already
# known, written as a literal.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a brand-new worksheet right after "总计" (i.e. before the
#    sheet currently in position 2, "2021-Q3") and name it "2022-Q3".
#    This pushes 2021-Q3 / 2021-Q2 / 2021-Q1 / 2020-Q4 one slot to the
#    right, matching the diff.
# ---------------------------------------------------------------------
$firstQuarterSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($firstQuarterSheet)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Populate the new sheet with the fund-holding table for 2022-Q3.
#    Layout mirrors the other quarterly sheets (2021-Q3, etc.) except
#    the "fund size" header reads "基金规模" for this quarter.
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

# "000270" must stay text (its leading zeros would be lost as a number);
# force text entry the same way Excel does via the Text number format,
# then drop back to the default style so the cell isn't left with a
# lingering custom format.
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "000270"
$newSheet.Range("B2").Style = "Normal"

$newSheet.Range("C2").Value = "建信灵活配置混合"

# D2/E2/F2/G2 are numeric-looking but must stay text too, matching the
# other quarterly sheets (fund size / position % / market value are all
# stored as plain text there, not numbers).
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.27"
$newSheet.Range("E2").Value = "94.21"
$newSheet.Range("F2").Value = "0.91"
$newSheet.Range("G2").Value = "0.0207"
$newSheet.Range("D2:G2").Style = "Normal"

$newSheet.Range("H2").Value = 6

# Match the header / index-column formatting used on the other quarterly
# sheets (bold, centered, thin border - cell style index 2 in the
# original workbook) by copying it from an already-styled sheet.
$sourceStyle = $wb.Worksheets.Item(3)
$sourceStyle.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$sourceStyle.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Update the "总计" overview sheet: shift the existing rows down by
#    one and insert the new 2022-Q3 summary row at the top of the data
#    (row 2), renumbering the index column (A) 0..4 along the way.
#
#    Original (rows 2-5): 2021-Q3/2/0.06, 2021-Q2/5/0.32,
#                          2021-Q1/6/0.47, 2020-Q4/3/0.18
#    New (rows 2-6):       2022-Q3/1/0.02, 2021-Q3/2/0.06,
#                          2021-Q2/5/0.32, 2021-Q1/6/0.47, 2020-Q4/3/0.18
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Give the new bottom row (row 6) the same index-column formatting as
# the other "总计" index cells before writing into it.
$totalSheet.Cells.Item(2, 1).Copy()
$totalSheet.Cells.Item(6, 1).PasteSpecial(-4122)

$dates  = @("2022-Q3", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$counts = @(1, 2, 5, 6, 3)
$values = @(0.02, 0.06, 0.32, 0.47, 0.18)

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $dates[$i]
    $totalSheet.Cells.Item($r, 3).Value = $counts[$i]
    $totalSheet.Cells.Item($r, 4).Value = $values[$i]
}
